# ADD results from server
# Update cached calculation results (investment cost values) on each year sheet.

$wb = $excel.ActiveWorkbook

$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("B2").Value = 9739.537847600008
$ws2025.Range("E2").Value = 289823.7596598056
$ws2025.Range("I2").Value = 161752.8135478
$ws2025.Range("L2").Value = 485245.29503538
$ws2025.Range("M2").Value = 105905.87968015
$ws2025.Range("N2").Value = 70831.955579581
$ws2025.Range("O2").Value = 69610.4422391004

$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("B2").Value = 47386.06393082884
$ws2030.Range("E2").Value = 271236.7992183856
$ws2030.Range("I2").Value = 280426.171173861
$ws2030.Range("L2").Value = 184420.4799505123
$ws2030.Range("M2").Value = 113936.92264746
$ws2030.Range("N2").Value = 33931.8246116005
$ws2030.Range("O2").Value = 50485.47232467777

$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("A2").Value = 28619.61401238371
$ws2035.Range("B2").Value = 23143.29485244409
$ws2035.Range("E2").Value = 111916.8406725409
$ws2035.Range("I2").Value = 150385.2728707001
$ws2035.Range("M2").Value = 34803.41203795493
$ws2035.Range("N2").Value = 44938.11408779013
$ws2035.Range("O2").Value = 26938.31306104351

$ws2040 = $wb.Worksheets.Item("2040")
$ws2040.Range("N2").Value = 1014.766490779938

$ws2045 = $wb.Worksheets.Item("2045")
$ws2045.Range("A2").Value = 34409.11717595647
$ws2045.Range("N2").Value = 5182.698656944208
$ws2045.Range("O2").Value = 22972.54525065906
